# MON-551 add nested joins
# Adds a new column H to the joinable_dummies worksheet containing the
# nested-join header/values ("another_joinable_relations.nested_joinable_relations.foreign_field")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new nested-join column
$ws.Range("H1").Value = "another_joinable_relations.nested_joinable_relations.foreign_field"

# Data rows mirror the joinable_relations.foreign_field values
$ws.Range("H2").Value = "foreign_field 1"
$ws.Range("H3").Value = "foreign_field 2"
$ws.Range("H4").Value = "foreign_field 3"

# Update the active selection to reflect where editing left off
$ws.Range("F11").Select()
